{"js": "// Remove the \"1097178 - Jo\u00e3o Batista de Almeida e Silva\" bullet paragraph\n// that used to sit right after the \"Docente(s) Respons\u00e1vel(eis)\" heading.\nconst body = context.document.body;\n\nconst results = body.search(\"1097178 - Jo\u00e3o Batista de Almeida e Silva\", {\n  matchCase: true,\n  matchWholeWord: false\n});\nresults.load(\"items\");\nawait context.sync();\n\nfor (let i = 0; i < results.items.length; i++) {\n  const hitParagraphs = results.items[i].paragraphs;\n  hitParagraphs.load(\"items\");\n  await context.sync();\n  for (let j = 0; j < hitParagraphs.items.length; j++) {\n    hitParagraphs.items[j].delete();\n  }\n}\n\nawait context.sync();\n", "ps1": "# Remove the \"1097178 - Jo\u00e3o Batista de Almeida e Silva\" bullet paragraph\n# that used to sit right after the \"Docente(s) Respons\u00e1vel(eis)\" heading.\n$d = $word.ActiveDocument\n\n$target = \"1097178 - Jo\u00e3o Batista de Almeida e Silva\"\n$wdParagraph = 4\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Text = $target\n$find.MatchCase = $true\n$find.Forward = $true\n\nwhile ($find.Execute()) {\n    $hit = $find.Parent\n    $hit.Expand($wdParagraph) | Out-Null\n    $hit.Delete()\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Text = $target\n    $find.MatchCase = $true\n    $find.Forward = $true\n}\n"}
